# Fixed bug with shop not correctly counting unaccounted-for stock:
# deduct the shopping-list quantities for "Bill" from the Inventory
# stock counts, and record Bill's order on the Shopping List / Expenses /
# Items Not Found sheets.

$wb = $excel.ActiveWorkbook

# --- Inventory sheet: correct Stock (column D) quantities ---------------
$invSheet = $wb.Worksheets.Item("Inventory")

$invSheet.Range("D2").Value  = 989
$invSheet.Range("D3").Value  = 701
$invSheet.Range("D5").Value  = 740
$invSheet.Range("D6").Value  = 9974
$invSheet.Range("D7").Value  = 9790
$invSheet.Range("D9").Value  = 85
$invSheet.Range("D10").Value = 87
$invSheet.Range("D11").Value = 5

$invSheet.Range("E18").Select()

# --- Shopping List sheet: add Bill's Laptops order -----------------------
$shopSheet = $wb.Worksheets.Item("Shopping List")

$shopSheet.Cells.Item(12, 1).Value = "Bill"
$shopSheet.Cells.Item(12, 2).Value = "Laptops"
$shopSheet.Cells.Item(12, 3).Value = 30

$shopSheet.Range("D16").Select()

# --- Expenses sheet: add Bill's total -------------------------------------
$expSheet = $wb.Worksheets.Item("Expenses")

$expSheet.Cells.Item(5, 1).Value = "Bill"
$expSheet.Cells.Item(5, 2).Value = 3400

$expSheet.Range("B7").Select()

# --- Items Not Found sheet: add Bill's unaccounted Laptops quantity ------
$notFoundSheet = $wb.Worksheets.Item("Items Not Found")

$notFoundSheet.Cells.Item(3, 1).Value = "Bill"
$notFoundSheet.Cells.Item(3, 2).Value = "Laptops"
$notFoundSheet.Cells.Item(3, 3).Value = 25

# Re-activate the Inventory sheet/tab as the visible sheet, matching the
# original workbook's tabSelected state.
$invSheet.Activate()
